$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-03-05 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-06 Thursday", 2)
$d.Content.Find.Execute("390÷2=195, 0", $true, $false, $false, $false, $false, $true, 1, $false, "559÷7=79, 6", 2)
$d.Content.Find.Execute("505÷8=63, 1", $true, $false, $false, $false, $false, $true, 1, $false, "401÷9=44, 5", 2)
$d.Content.Find.Execute("675÷7=96, 3", $true, $false, $false, $false, $false, $true, 1, $false, "457÷9=50, 7", 2)
$d.Content.Find.Execute("835÷9=92, 7", $true, $false, $false, $false, $false, $true, 1, $false, "377÷5=75, 2", 2)
$d.Content.Find.Execute("845÷7=120, 5", $true, $false, $false, $false, $false, $true, 1, $false, "141÷5=28, 1", 2)
$d.Content.Find.Execute("835÷5=167, 0", $true, $false, $false, $false, $false, $true, 1, $false, "196÷9=21, 7", 2)
$d.Content.Find.Execute("670÷6=111, 4", $true, $false, $false, $false, $false, $true, 1, $false, "821÷8=102, 5", 2)
$d.Content.Find.Execute("678÷5=135, 3", $true, $false, $false, $false, $false, $true, 1, $false, "889÷9=98, 7", 2)
$d.Content.Find.Execute("727÷8=90, 7", $true, $false, $false, $false, $false, $true, 1, $false, "983÷3=327, 2", 2)
$d.Content.Find.Execute("531÷9=59, 0", $true, $false, $false, $false, $false, $true, 1, $false, "725÷5=145, 0", 2)
$d.Content.Find.Execute("789÷7=112, 5", $true, $false, $false, $false, $false, $true, 1, $false, "141÷5=28, 1", 2)
$d.Content.Find.Execute("267÷4=66, 3", $true, $false, $false, $false, $false, $true, 1, $false, "259÷4=64, 3", 2)
$d.Content.Find.Execute("786÷8=98, 2", $true, $false, $false, $false, $false, $true, 1, $false, "464÷9=51, 5", 2)
$d.Content.Find.Execute("177÷9=19, 6", $true, $false, $false, $false, $false, $true, 1, $false, "964÷2=482, 0", 2)
$d.Content.Find.Execute("840÷2=420, 0", $true, $false, $false, $false, $false, $true, 1, $false, "532÷6=88, 4", 2)
$d.Content.Find.Execute("180÷7=25, 5", $true, $false, $false, $false, $false, $true, 1, $false, "514÷7=73, 3", 2)
$d.Content.Find.Execute("135÷8=16, 7", $true, $false, $false, $false, $false, $true, 1, $false, "854÷2=427, 0", 2)
$d.Content.Find.Execute("945÷7=135, 0", $true, $false, $false, $false, $false, $true, 1, $false, "669÷8=83, 5", 2)
$d.Content.Find.Execute("434÷2=217, 0", $true, $false, $false, $false, $false, $true, 1, $false, "130÷5=26, 0", 2)
$d.Content.Find.Execute("941÷7=134, 3", $true, $false, $false, $false, $false, $true, 1, $false, "501÷9=55, 6", 2)
$d.Content.Find.Execute("414÷6=69, 0", $true, $false, $false, $false, $false, $true, 1, $false, "766÷2=383, 0", 2)
$d.Content.Find.Execute("795÷5=159, 0", $true, $false, $false, $false, $false, $true, 1, $false, "119÷9=13, 2", 2)
$d.Content.Find.Execute("349÷5=69, 4", $true, $false, $false, $false, $false, $true, 1, $false, "893÷7=127, 4", 2)
$d.Content.Find.Execute("867÷4=216, 3", $true, $false, $false, $false, $false, $true, 1, $false, "809÷6=134, 5", 2)
$d.Content.Find.Execute("163÷9=18, 1", $true, $false, $false, $false, $false, $true, 1, $false, "639÷8=79, 7", 2)
